$d = $word.ActiveDocument

$pairs = @(
    @('2023-10-19 Thursday', '2023-10-20 Friday'),
    @('10+45=', '39+59='),
    @('43-16=', '81-66='),
    @('76-59=', '93-34='),
    @('57-13=', '13+28='),
    @('5-4=', '41-4='),
    @('85+4=', '92-90='),
    @('24+57=', '60+21='),
    @('97-75=', '47-20='),
    @('80-27=', '76-73='),
    @('25+1=', '93-65='),
    @('89-1=', '24-11='),
    @('96-25=', '27-21='),
    @('61+26=', '72-54='),
    @('91-55=', '81+9='),
    @('61+14=', '87-11='),
    @('10+55=', '23-21='),
    @('63+2=', '93-76='),
    @('9-0=', '85-28='),
    @('23+21=', '25+37='),
    @('87-28=', '36+50='),
    @('93-5=', '65+19='),
    @('1+59=', '32+35='),
    @('97-45=', '70+26='),
    @('11+41=', '67-28='),
    @('54+23=', '7+41='),
    @('67+12=', '91-48='),
    @('29-14=', '30+0='),
    @('24+9=', '12-3='),
    @('88-32=', '39+5='),
    @('38+25=', '6+60='),
    @('55+41=', '15+19='),
    @('4-0=', '77-0='),
    @('32+57=', '79-78='),
    @('38+27=', '95-4='),
    @('16-4=', '82-58='),
    @('93-63=', '72+15='),
    @('68+0=', '22-3='),
    @('76-72=', '82-31='),
    @('55+37=', '95-81='),
    @('40+1=', '56-44='),
    @('8+66=', '24+13='),
    @('41-36=', '60+30='),
    @('77-48=', '19+18='),
    @('25+43=', '45-2='),
    @('82-6=', '96-1='),
    @('82+17=', '99-22='),
    @('45-24=', '56+42='),
    @('16-2=', '72-66='),
    @('92-50=', '55-42='),
    @('30-1=', '82+6='),
    @('29+28=', '12+4='),
    @('6+5=', '15+80='),
    @('27-4=', '68-60='),
    @('58-4=', '90-34='),
    @('40+56=', '29-26='),
    @('43-25=', '18+10='),
    @('5-0=', '7+38='),
    @('71-64=', '12+21='),
    @('64-48=', '27+3='),
    @('23-7=', '78-62='),
    @('26+70=', '0+13='),
    @('3+17=', '18+29='),
    @('49-19=', '89-67='),
    @('55+32=', '58-31='),
    @('36+63=', '72-24='),
    @('98-73=', '69-53='),
    @('61-18=', '15+44='),
    @('13+38=', '30-29='),
    @('77-25=', '21-20='),
    @('83-55=', '59-32='),
    @('33+10=', '31-9='),
    @('27+6=', '93-4='),
    @('77-33=', '97-29='),
    @('13+1=', '0+49='),
    @('46+29=', '81-36='),
    @('46-46=', '85-85='),
    @('70+20=', '60-5='),
    @('68+10=', '42+49='),
    @('71-40=', '88-85='),
    @('74+16=', '70-51='),
    @('37-26=', '62+32='),
    @('51-2=', '13-2='),
    @('7+88=', '42-14='),
    @('18+25=', '42+27='),
    @('62+24=', '96-2='),
    @('22+36=', '28+43='),
    @('88-12=', '58+41='),
    @('39+11=', '48+12='),
    @('34+50=', '96-67='),
    @('46-45=', '3-3='),
    @('76-6=', '53-10='),
    @('77-23=', '39+42='),
    @('44+3=', '48-9='),
    @('7+55=', '32-17='),
    @('90-74=', '10+58='),
    @('91-57=', '0+47='),
    @('9+46=', '27+67='),
    @('41-10=', '52+6='),
    @('55-17=', '59-22='),
    @('87-10=', '72+4='),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()